$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.267.02"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "1.657.70"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("E4").Value = "  -0.72%  "
$ws.Range("D5").Value = "219.46"
$ws.Range("E5").Value = "  -0.79%  "
$ws.Range("D6").Value = "0.5246"
$ws.Range("E6").Value = "  -1.61%  "
$ws.Range("E7").Value = "  -0.69%  "
$ws.Range("E8").Value = "  +0.37%  "
$ws.Range("D9").Value = "0.06373"
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").Value = "20.74"
$ws.Range("E10").Value = "  -0.55%  "
$ws.Range("D11").Value = "0.07719"
$ws.Range("E11").Value = "  -1.84%  "
$ws.Range("D12").Value = "4.603"
$ws.Range("E12").Value = "  +1.66%  "
$ws.Range("D13").Value = "1.592.51"
$ws.Range("E13").Value = "  -4.88%  "
$ws.Range("D14").Value = "1.885.90"
$ws.Range("E14").Value = "  -0.80%  "
$ws.Range("E15").Value = "  +0.84%  "
$ws.Range("D16").Value = "0.0₅8269"
$ws.Range("E16").Value = "  +0.99%  "
$ws.Range("D17").Value = "65.54"
$ws.Range("E17").Value = "  -0.93%  "
$ws.Range("D18").Value = "26.265.12"
$ws.Range("E18").Value = "  -0.60%  "
$ws.Range("E19").Value = "  -0.66%  "
$ws.Range("D20").Value = "4.695"
$ws.Range("E20").Value = "  -0.51%  "
$ws.Range("D21").Value = "193.31"
$ws.Range("E21").Value = "  -2.33%  "
$ws.Range("D22").Value = "10.45"
$ws.Range("E22").Value = "  +1.46%  "
$ws.Range("E23").Value = "  -1.02%  "
$ws.Range("E24").Value = "  -0.67%  "
$ws.Range("D25").Value = "143.29"
$ws.Range("E25").Value = "  -1.71%  "
$ws.Range("D26").Value = "0.1204"
$ws.Range("E26").Value = "  -1.76%  "
$ws.Range("D27").Value = "7.302"
$ws.Range("E27").Value = "  +0.67%  "
$ws.Range("D28").Value = "15.95"
$ws.Range("E28").Value = "  -1.84%  "
$ws.Range("D29").Value = "1.513"
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("D30").Value = "0.05664"
$ws.Range("E30").Value = "  -4.28%  "
$ws.Range("D31").Value = "1.278"
$ws.Range("E31").Value = "  -0.84%  "
$ws.Range("D32").Value = "3.515"
$ws.Range("E32").Value = "  -1.37%  "
$ws.Range("D33").Value = "3.359"
$ws.Range("E33").Value = "  +0.79%  "
$ws.Range("D34").Value = "1.587"
$ws.Range("E34").Value = "  -1.40%  "
$ws.Range("D35").Value = "2.804"
$ws.Range("E35").Value = "  -1.20%  "
$ws.Range("D36").Value = "0.9499"
$ws.Range("E36").Value = "  -1.96%  "
$ws.Range("D37").Value = "2.413"
$ws.Range("E37").Value = "  -1.06%  "
$ws.Range("D38").Value = "0.5778"
$ws.Range("E38").Value = "  -1.04%  "
$ws.Range("E39").Value = "  -0.85%  "
$ws.Range("D40").Value = "5.976"
$ws.Range("E40").Value = "  +0.66%  "
$ws.Range("E41").Value = "  -1.18%  "
$ws.Range("D42").Value = "0.8473"
$ws.Range("E42").Value = "  -2.03%  "
$ws.Range("E43").Value = "  -0.75%  "
$ws.Range("D44").Value = "101.97"
$ws.Range("E44").Value = "  -1.10%  "
$ws.Range("D45").Value = "1.019.40"
$ws.Range("E45").Value = "  -5.61%  "
$ws.Range("D46").Value = "1.796.91"
$ws.Range("E46").Value = "  -0.80%  "
$ws.Range("D47").Value = "58.35"
$ws.Range("E47").Value = "  -0.26%  "
$ws.Range("D48").Value = "0.0₈107"
$ws.Range("E48").Value = "  -0.15%  "
$ws.Range("E49").Value = "  -1.12%  "
$ws.Range("D50").Value = "0.05329"
$ws.Range("E50").Value = "  +3.26%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "8.024"
$ws.Range("E51").Value = "  -0.08%  "
